# Applies the change described by the diff:
#  - "Test Cases" sheet: Runmode for TC_IncidentReportSaveSubmit (row 6) changes Yes -> No
#  - "Test Cases" sheet: four new test-case rows are appended (23-26) for the new
#    investigation / acknowledge / turn-around-time test cases
#  - a new worksheet "TC_IR_TAT" is added after "TC_QualityDashboardReAssign"
#  - view/selection state is updated to reflect the edits made on "Test Cases"

$wb = $excel.ActiveWorkbook

$testCases = $wb.Worksheets.Item("Test Cases")

# Runmode for TC_IncidentReportSaveSubmit flips from Yes to No because the newly
# added TC_IR_TAT case takes over as the active ("Yes") run.
$testCases.Range("C6").Value = "No"

# New rows are authored in this exact order (matches how the shared string table
# ends up populated): A24, A23, A26, B23, B24, B25, B26, A25, then the Runmode column.
$testCases.Range("A24").Value = "TC_IRinvestigationRecordInputs"
$testCases.Range("A23").Value = "TC_IRinvestigation"
$testCases.Range("A26").Value = "TC_IR_TAT"
$testCases.Range("B23").Value = "Tests the investigation functionality with asking details"
$testCases.Range("B24").Value = "Tests the investigation functionality with filling details"
$testCases.Range("B25").Value = "Tests the acknowledge functionality of the incidents"
$testCases.Range("B26").Value = "Tests the Turn Around Time functionality of the incidents"
$testCases.Range("A25").Value = "TC_IRAcknowledge"

$testCases.Range("C23").Value = "No"
$testCases.Range("C24").Value = "No"
$testCases.Range("C25").Value = "No"
$testCases.Range("C26").Value = "Yes"

# New worksheet "TC_IR_TAT" added after the last existing sheet (TC_QualityDashboardReAssign)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "TC_IR_TAT"

# Reuse formatting from the analogous header/value cells on TC_QualityDashboardReAssign
# so the same fill/number-format styles are shared, then center the text.
$reassignSheet = $wb.Worksheets.Item("TC_QualityDashboardReAssign")
$reassignSheet.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null
$reassignSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122) | Out-Null

$newSheet.Range("A1").Value = "IR Code"
$newSheet.Range("A2").Value = "No"

$newSheet.Range("A1").HorizontalAlignment = -4108
$newSheet.Range("A2").HorizontalAlignment = -4108

$newSheet.Columns.Item(1).ColumnWidth = 24.5

$newSheet.Range("B2").Select() | Out-Null

# Make "Test Cases" the active/selected sheet again, scrolled/selected on the newly
# added rows.
$testCases.Activate()
$testCases.Range("A26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8

Write-Host "edit applied"
